$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.268.07"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "2.243.93"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.17"
$ws.Range("E5").Value = "  +2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  -5.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.47"
$ws.Range("E7").Value = "  -4.46%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  -6.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.32"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "36.68"
$ws.Range("E12").Value = "  +11.25%  "
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("E14").Value = "  -6.72%  "
$ws.Range("D15").Value = "2.578.19"
$ws.Range("E15").Value = "  -4.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.04"
$ws.Range("E16").Value = "  -6.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.868"
$ws.Range("E17").Value = "  -3.69%  "
$ws.Range("D18").Value = "2.242.42"
$ws.Range("E18").Value = "  -4.31%  "
$ws.Range("D19").Value = "42.153.33"
$ws.Range("E19").Value = "  -3.68%  "
$ws.Range("D20").Value = "0.0₃0973"
$ws.Range("E20").Value = "  -4.81%  "
$ws.Range("E21").Value = "  -5.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.26"
$ws.Range("E22").Value = "  -6.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.42"
$ws.Range("E23").Value = "  -5.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.05"
$ws.Range("E24").Value = "  +10.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.64"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.48"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "172.39"
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.60"
$ws.Range("E31").Value = "  -7.30%  "
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -5.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0718"
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.31"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.71"
$ws.Range("E36").Value = "  -7.31%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.00"
$ws.Range("E38").Value = "  +22.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0284"
$ws.Range("E39").Value = "  +4.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.31"
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "67.74"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.91"
$ws.Range("E42").Value = "  -7.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.31"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.05"
$ws.Range("E44").Value = "  -9.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.102"
$ws.Range("E45").Value = "  -3.51%  "
$ws.Range("B46").Value = "SynthetixNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.67"
$ws.Range("E46").Value = "  +12.82%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.191"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.20"
$ws.Range("E49").Value = "  -1.60%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.10"
$ws.Range("E50").Value = "  +7.86%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.34"
$ws.Range("E51").Value = "  -2.89%  "
